# Scheduled-runner market data refresh: updates currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H:N) across all 8 job sheets with
# freshly-polled values. Plain static values -- no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1394
$ws.Range("I2").Value = 1394
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1394
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1281

$ws.Range("H17").Value = 1909
$ws.Range("J17").Value = 2833.1667
$ws.Range("L17").Value = 8499.500100000001
$ws.Range("N17").Value = -8835.500100000001

$ws.Range("H43").Value = 1132.6666
$ws.Range("I43").Value = 1099.5
$ws.Range("K43").Value = 1099.5
$ws.Range("M43").Value = -1030.5

$ws.Range("H92").Value = 347
$ws.Range("I92").Value = 393.5625
$ws.Range("J92").Value = 240.57143
$ws.Range("K92").Value = 393.5625
$ws.Range("L92").Value = 240.57143
$ws.Range("M92").Value = 854.4375
$ws.Range("N92").Value = -2736.57143

$ws.Range("H98").Value = 2673.5833
$ws.Range("I98").Value = 1077.5
$ws.Range("J98").Value = 2992.8
$ws.Range("K98").Value = 1077.5
$ws.Range("L98").Value = 2992.8
$ws.Range("M98").Value = 420.5
$ws.Range("N98").Value = -5988.8

$ws.Range("H122").Value = 2673.5833
$ws.Range("I122").Value = 1077.5
$ws.Range("J122").Value = 2992.8
$ws.Range("K122").Value = 3232.5
$ws.Range("L122").Value = 8978.400000000001
$ws.Range("M122").Value = -782.5
$ws.Range("N122").Value = -13878.4

$ws.Range("H134").Value = 70655.60000000001
$ws.Range("J134").Value = 70655.60000000001
$ws.Range("L134").Value = 70655.60000000001
$ws.Range("N134").Value = -80795.60000000001

$ws.Range("H135").Value = 1854.6842
$ws.Range("I135").Value = 1452.4375
$ws.Range("K135").Value = 13071.9375
$ws.Range("M135").Value = -10536.9375

$ws.Range("H137").Value = 455831.44
$ws.Range("J137").Value = 1321248.9
$ws.Range("L137").Value = 3963746.7
$ws.Range("N137").Value = -3968846.7

$ws.Range("H138").Value = 2736.225
$ws.Range("I138").Value = 2136.6
$ws.Range("J138").Value = 3335.85
$ws.Range("K138").Value = 6409.799999999999
$ws.Range("L138").Value = 10007.55
$ws.Range("M138").Value = -1269.799999999999
$ws.Range("N138").Value = -20287.55

$ws.Range("H141").Value = 2132.25
$ws.Range("I141").Value = 2316.389
$ws.Range("J141").Value = 475
$ws.Range("K141").Value = 6949.167
$ws.Range("L141").Value = 1425
$ws.Range("M141").Value = -1769.167
$ws.Range("N141").Value = -11785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6001.6616
$ws.Range("I32").Value = 2952.4211
$ws.Range("J32").Value = 21802.273
$ws.Range("K32").Value = 2952.4211
$ws.Range("L32").Value = 21802.273
$ws.Range("M32").Value = -2665.4211
$ws.Range("N32").Value = -22376.273

$ws.Range("H97").Value = 661.86664
$ws.Range("I97").Value = 623.5833
$ws.Range("K97").Value = 623.5833
$ws.Range("M97").Value = -127.5833

$ws.Range("H110").Value = 1667.8889
$ws.Range("I110").Value = 1015.8571
$ws.Range("K110").Value = 1015.8571
$ws.Range("M110").Value = 1029.1429

$ws.Range("H132").Value = 1801.0588
$ws.Range("I132").Value = 1339.1666
$ws.Range("J132").Value = 2909.6
$ws.Range("K132").Value = 4017.4998
$ws.Range("L132").Value = 8728.799999999999
$ws.Range("M132").Value = -1487.4998
$ws.Range("N132").Value = -13788.8

$ws.Range("H134").Value = 116969.336
$ws.Range("J134").Value = 116969.336
$ws.Range("L134").Value = 116969.336
$ws.Range("N134").Value = -127109.336

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 698.625
$ws.Range("I5").Value = 631
$ws.Range("K5").Value = 631
$ws.Range("M5").Value = -518

$ws.Range("H94").Value = 1492.28
$ws.Range("I94").Value = 1262.375
$ws.Range("K94").Value = 1262.375
$ws.Range("M94").Value = -811.375

$ws.Range("H102").Value = 17499.5
$ws.Range("I102").Value = 15000
$ws.Range("J102").Value = 19999
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 19999
$ws.Range("M102").Value = -11755
$ws.Range("N102").Value = -26489

$ws.Range("H107").Value = 1447.0667
$ws.Range("I107").Value = 1246.7273
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 1246.7273
$ws.Range("L107").Value = 1998
$ws.Range("M107").Value = 673.2727
$ws.Range("N107").Value = -5838

$ws.Range("H134").Value = 2605.6086
$ws.Range("I134").Value = 1462
$ws.Range("J134").Value = 4749.875
$ws.Range("K134").Value = 4386
$ws.Range("L134").Value = 14249.625
$ws.Range("M134").Value = -1851
$ws.Range("N134").Value = -19319.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6263.5454
$ws.Range("I7").Value = 8473.083000000001
$ws.Range("J7").Value = 5000.952
$ws.Range("K7").Value = 8473.083000000001
$ws.Range("L7").Value = 5000.952
$ws.Range("M7").Value = -8360.083000000001
$ws.Range("N7").Value = -5226.952

$ws.Range("H12").Value = 785.4286
$ws.Range("I12").Value = 583.1667
$ws.Range("J12").Value = 1999
$ws.Range("K12").Value = 583.1667
$ws.Range("L12").Value = 1999
$ws.Range("M12").Value = -413.1667
$ws.Range("N12").Value = -2339

$ws.Range("H31").Value = 2026.238
$ws.Range("I31").Value = 1709.3889
$ws.Range("K31").Value = 1709.3889
$ws.Range("M31").Value = -1414.3889

$ws.Range("H34").Value = 2026.238
$ws.Range("I34").Value = 1709.3889
$ws.Range("K34").Value = 1709.3889
$ws.Range("M34").Value = -1507.3889

$ws.Range("H105").Value = 38656.266
$ws.Range("I105").Value = 59468.633
$ws.Range("J105").Value = 2707.6365
$ws.Range("K105").Value = 59468.633
$ws.Range("L105").Value = 2707.6365
$ws.Range("M105").Value = -57721.633
$ws.Range("N105").Value = -6201.636500000001

$ws.Range("H132").Value = 2390.2173
$ws.Range("I132").Value = 2258.25
$ws.Range("J132").Value = 2691.8572
$ws.Range("K132").Value = 6774.75
$ws.Range("L132").Value = 8075.571599999999
$ws.Range("M132").Value = -4244.75
$ws.Range("N132").Value = -13135.5716

$ws.Range("H134").Value = 29334.842
$ws.Range("I134").Value = 2844.3872
$ws.Range("K134").Value = 8533.161599999999
$ws.Range("M134").Value = -5998.161599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 119816.63
$ws.Range("I8").Value = 119816.63
$ws.Range("K8").Value = 359449.89
$ws.Range("M8").Value = -359310.89

$ws.Range("H121").Value = 3031
$ws.Range("J121").Value = 3572.111
$ws.Range("L121").Value = 10716.333
$ws.Range("N121").Value = -13336.333

$ws.Range("H131").Value = 44579.39
$ws.Range("I131").Value = 72182.71000000001
$ws.Range("J131").Value = 1640.8889
$ws.Range("K131").Value = 216548.13
$ws.Range("L131").Value = 4922.6667
$ws.Range("M131").Value = -211508.13
$ws.Range("N131").Value = -15002.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 492.22223
$ws.Range("I2").Value = 740.1177
$ws.Range("J2").Value = 270.42105
$ws.Range("K2").Value = 740.1177
$ws.Range("L2").Value = 270.42105
$ws.Range("M2").Value = -627.1177
$ws.Range("N2").Value = -496.42105

$ws.Range("H97").Value = 1020.5161
$ws.Range("I97").Value = 559
$ws.Range("J97").Value = 4135.75
$ws.Range("K97").Value = 559
$ws.Range("L97").Value = 4135.75
$ws.Range("M97").Value = -63
$ws.Range("N97").Value = -5127.75

$ws.Range("H102").Value = 1642.1111
$ws.Range("I102").Value = 1680.3529
$ws.Range("J102").Value = 992
$ws.Range("K102").Value = 1680.3529
$ws.Range("L102").Value = 992
$ws.Range("M102").Value = -58.35290000000009
$ws.Range("N102").Value = -4236

$ws.Range("H122").Value = 144160.47
$ws.Range("I122").Value = 202773.25
$ws.Range("J122").Value = 3489.8
$ws.Range("K122").Value = 608319.75
$ws.Range("L122").Value = 10469.4
$ws.Range("M122").Value = -605869.75
$ws.Range("N122").Value = -15369.4

$ws.Range("H132").Value = 8005
$ws.Range("J132").Value = 7439.9
$ws.Range("L132").Value = 22319.7
$ws.Range("N132").Value = -27379.7

$ws.Range("H134").Value = 63499.5
$ws.Range("J134").Value = 63499.5
$ws.Range("L134").Value = 190498.5
$ws.Range("N134").Value = -195568.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 687
$ws.Range("I16").Value = 622.7619
$ws.Range("K16").Value = 622.7619
$ws.Range("M16").Value = -452.7619

$ws.Range("H55").Value = 4416.3516
$ws.Range("I55").Value = 1002
$ws.Range("J55").Value = 15038.777
$ws.Range("K55").Value = 1002
$ws.Range("L55").Value = 15038.777
$ws.Range("M55").Value = -829
$ws.Range("N55").Value = -15384.777

$ws.Range("H132").Value = 2261.5264
$ws.Range("I132").Value = 2027.5883
$ws.Range("K132").Value = 6082.7649
$ws.Range("M132").Value = -3552.7649

$ws.Range("H136").Value = 1827.6666
$ws.Range("I136").Value = 1589.6111
$ws.Range("J136").Value = 2303.7778
$ws.Range("K136").Value = 4768.8333
$ws.Range("L136").Value = 6911.3334
$ws.Range("M136").Value = -2218.8333
$ws.Range("N136").Value = -12011.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1116231
$ws.Range("I132").Value = 1251.7667
$ws.Range("J132").Value = 4832828.5
$ws.Range("K132").Value = 3755.300099999999
$ws.Range("L132").Value = 14498485.5
$ws.Range("M132").Value = -1225.300099999999
$ws.Range("N132").Value = -14503545.5

$ws.Range("H136").Value = 1947.2858
$ws.Range("J136").Value = 3249.1667
$ws.Range("L136").Value = 9747.500100000001
$ws.Range("N136").Value = -14847.5001
